# Applies the "Added new evpwd params" edit to the workbook.
# Target worksheet: "evp-wd_f" (4th sheet in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("evp-wd_f")

# --- Row 7: correct / update an existing data row ---
$ws.Range("C7").Value = "8; 7"
$ws.Range("D7").Value = 4.1862000000000004
$ws.Range("E7").Value = 84.548000000000002
$ws.Range("F7").Value = 2.1122999999999998
$ws.Range("G7").Value = 4.7751999999999999
$ws.Range("H7").Value = 1574.3
$ws.Range("I7").Value = 1.6472
$ws.Range("J7").Value = 0.34209000000000001
$ws.Range("K7").Value = 3.5358000000000001
$ws.Range("L7").Value = 10.343
$ws.Range("M7").Value = 0.55237999999999998
$ws.Range("N7").Value = 2.4053

# Row 7 (D:N) gets re-centred horizontally only (no vertical centring, no border)
$rowSeven = $ws.Range("D7:N7")
$rowSeven.HorizontalAlignment = -4108
$rowSeven.VerticalAlignment = -4107
$rowSeven.Borders.LineStyle = -4142

# --- Row 34: fill in new evpwd params ---
$ws.Range("C34").Value = "7; 5"
$ws.Range("I34").Value = 1.6442000000000001
$ws.Range("J34").Value = 0.55276000000000003
$ws.Range("K34").Value = 4.5461
$ws.Range("L34").Value = 8.8484999999999996
$ws.Range("M34").Value = 0.20913000000000001
$ws.Range("N34").Value = 2.6558000000000002

# --- Row 11: fill in new evpwd params ---
$ws.Range("C11").Value = "5; 6"
$ws.Range("I11").Value = 2.2593000000000001
$ws.Range("J11").Value = 0.63022
$ws.Range("K11").Value = 4.9957000000000003
$ws.Range("L11").Value = 1.0884
$ws.Range("M11").Value = 0.051331000000000002
$ws.Range("N11").Value = 3.407

# --- Row 23: fill in new evpwd params ---
$ws.Range("C23").Value = "4; 7"
$ws.Range("I23").Value = 1.5074000000000001
$ws.Range("J23").Value = 0.46747
$ws.Range("K23").Value = 4.2224000000000004
$ws.Range("L23").Value = 7.0061
$ws.Range("M23").Value = 0.72074000000000005
$ws.Range("N23").Value = 4.0982000000000003

# --- Row 10: fill in new evpwd params ---
$ws.Range("C10").Value = "7; 7"
$ws.Range("I10").Value = 1.1933
$ws.Range("J10").Value = 0.18759999999999999
$ws.Range("K10").Value = 2.4325999999999999
$ws.Range("L10").Value = 2.4598
$ws.Range("M10").Value = 0.27822999999999998
$ws.Range("N10").Value = 3.0878999999999999

# --- Row 22: fill in new evpwd params ---
$ws.Range("C22").Value = "6; 7"
$ws.Range("I22").Value = 1.5924
$ws.Range("J22").Value = 0.33024999999999999
$ws.Range("K22").Value = 3.5438999999999998
$ws.Range("L22").Value = 2.0625
$ws.Range("M22").Value = 0.012094000000000001
$ws.Range("N22").Value = 2.7164000000000001

# --- Row 35: fill in new evpwd params ---
$ws.Range("C35").Value = "6; 7"
$ws.Range("I35").Value = 3.2982
$ws.Range("J35").Value = 0.69935000000000003
$ws.Range("K35").Value = 5.0255000000000001
$ws.Range("L35").Value = 2.2077
$ws.Range("M35").Value = 0.21060999999999999
$ws.Range("N35").Value = 3.0865

# --- Row 46: fill in new evpwd params ---
$ws.Range("C46").Value = "6; 10"
$ws.Range("I46").Value = 2.1859999999999999
$ws.Range("J46").Value = 0.20261000000000001
$ws.Range("K46").Value = 2.2957999999999998
$ws.Range("L46").Value = 8.5665999999999993
$ws.Range("M46").Value = 0.39078000000000002
$ws.Range("N46").Value = 2.7238000000000002

# --- Row 47: fill in new evpwd params ---
$ws.Range("C47").Value = "7; 9"
$ws.Range("I47").Value = 1.8718999999999999
$ws.Range("J47").Value = 0.38352999999999998
$ws.Range("K47").Value = 3.4464999999999999
$ws.Range("L47").Value = 3.8252999999999999
$ws.Range("M47").Value = 0.19616
$ws.Range("N47").Value = 2.4363999999999999

# Update the view state (scroll position / selected cell) on the active sheet
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("C48").Select()
